$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-point the three "Table_0" tables (slides 14, 15, 16) to the standard
#    built-in table style {F1635EA7-F42D-418E-89CB-86473E9AD9F4} instead of
#    the custom style {3B4239C7-1C0A-4D3D-A871-913F350EDA3A}.
# ---------------------------------------------------------------------------
$newStyleId = "{F1635EA7-F42D-418E-89CB-86473E9AD9F4}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's colour theme from "Integral / Red Violet" to
#    "Office Theme / Office" (the colours that used to live on the Notes
#    Master's theme now become the colours used by the slides' theme).
#    ThemeColorScheme.Colors(i).RGB uses the VBA RGB() packing (0xBBGGRR),
#    so each OOXML "RRGGBB" value below is byte-reversed before assignment.
# ---------------------------------------------------------------------------
function ToVbaRgb([string]$rrggbb) {
    $r = $rrggbb.Substring(0, 2)
    $g = $rrggbb.Substring(2, 2)
    $b = $rrggbb.Substring(4, 2)
    return [Convert]::ToInt32("$b$g$r", 16)
}

# Index order matches the OOXML <a:clrScheme> child order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToVbaRgb $officeColors[$i - 1]
}
